{"js": "// Remove the trailing \"Ver no Jupiter...\" and copyright paragraphs (and the\n// blank paragraph that separated them from the \"Requisitos\" section), while\n// leaving the final blank paragraph and page-break paragraph untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"LOM3049: ...\") that must stay.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"LOM3049: Termodin\u00e2mica de M\u00e1quinas (Requisito)\") {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error('Could not find anchor paragraph \"LOM3049: Termodin\u00e2mica de M\u00e1quinas (Requisito)\".');\n}\n\n// The three paragraphs directly following the anchor are expected to be:\n//   anchorIndex + 1 -> \"\" (blank separator)\n//   anchorIndex + 2 -> \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   anchorIndex + 3 -> \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\nconst expectedTexts = [\n  \"\",\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst toDelete = [];\nfor (let offset = 0; offset < expectedTexts.length; offset++) {\n  const i = anchorIndex + 1 + offset;\n  if (i < items.length && items[i].text === expectedTexts[offset]) {\n    toDelete.push(items[i]);\n  }\n}\n\n// Delete from the end backwards so earlier indices stay valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" and copyright paragraphs (and the\n# blank paragraph that separated them from the \"Requisitos\" section), while\n# leaving the final blank paragraph and page-break paragraph untouched.\n$d = $word.ActiveDocument\n\n$anchorText = \"LOM3049: Termodin\u00e2mica de M\u00e1quinas (Requisito)\"\n$expectedTexts = @(\n    \"\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph '$anchorText'.\"\n}\n\n# Collect the indices of the paragraphs right after the anchor that match the\n# expected texts (blank separator, \"Ver no Jupiter...\", \"\u00a9 2020...\").\n$toDelete = @()\nfor ($offset = 0; $offset -lt $expectedTexts.Count; $offset++) {\n    $i = $anchorIndex + 1 + $offset\n    if ($i -le $d.Paragraphs.Count) {\n        $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $expectedTexts[$offset]) {\n            $toDelete += $i\n        }\n    }\n}\n\n# Delete from the highest index down so earlier indices remain valid.\nfor ($j = $toDelete.Count - 1; $j -ge 0; $j--) {\n    $d.Paragraphs.Item($toDelete[$j]).Range.Delete()\n}\n"}
